# Updates cryptocurrency price/volume figures in the worksheet.
# Values are written as literal text (matching the sheet's existing text-typed
# cells), so numeric-looking prices are forced to Text via a temporary '@'
# number format that is reverted immediately after the write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text does not look like a plain number (percentages always
# carry '%' + padding spaces; some prices keep thousand-separator dots) can be
# assigned directly -- Excel keeps them as text automatically.
$plainUpdates = @{
    D2 = '72.989.95'
    E2 = '  +0.98%  '
    D3 = '4.060.80'
    E3 = '  +0.81%  '
    E4 = '  -0.10%  '
    E5 = '  +3.83%  '
    E6 = '  -0.58%  '
    D7 = '4.056.61'
    E7 = '  +0.97%  '
    E8 = '  -0.40%  '
    E9 = '  -0.10%  '
    E10 = '  +1.18%  '
    E11 = '  +0.37%  '
    E12 = '  +12.06%  '
    E13 = '  +1.04%  '
    E14 = '  +2.21%  '
    D15 = '4.704.47'
    E15 = '  +0.66%  '
    D16 = '4.055.44'
    E16 = '  +0.85%  '
    E17 = '  +2.83%  '
    E18 = '  +1.48%  '
    E19 = '  +1.43%  '
    E20 = '  -0.60%  '
    D21 = '72.841.25'
    E21 = '  +1.02%  '
    E22 = '  +3.78%  '
    E23 = '  -0.99%  '
    E24 = '  -0.22%  '
    E25 = '  +1.61%  '
    E26 = '  +1.90%  '
    E27 = '  +11.13%  '
    E28 = '  +2.03%  '
    E29 = '  +0.30%  '
    E30 = '  +1.52%  '
    E31 = '  +1.25%  '
    E32 = '  +13.97%  '
    E33 = '  +3.00%  '
    E34 = '  +1.53%  '
    E35 = '  +2.75%  '
    E36 = '  +14.97%  '
    E37 = '  +2.82%  '
    E38 = '  +5.65%  '
    E39 = '  +6.04%  '
    E40 = '  -3.00%  '
    E41 = '  -1.79%  '
    E42 = '  -1.21%  '
    E43 = '  +17.80%  '
    E44 = '  +0.00%  '
    E45 = '  +1.89%  '
    E46 = '  +0.01%  '
    E47 = '  +0.58%  '
    E48 = '  +3.11%  '
    E49 = '  +7.39%  '
    E50 = '  +3.96%  '
    E51 = '  -0.91%  '
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Cells whose new text parses as a plain number (e.g. "0.999", "559.81") need
# to be forced to Text so they don't silently become numeric cells; the sheet's
# original cells for these columns are plain text, not numbers.
$textForcedUpdates = @{
    D5 = '559.81'
    D6 = '152.20'
    D8 = '0.699'
    D9 = '0.999'
    D10 = '0.763'
    D12 = '53.66'
    D14 = '11.05'
    D17 = '14.60'
    D18 = '20.89'
    D22 = '449.53'
    D23 = '98.39'
    D24 = '3.56'
    D25 = '4.39'
    D26 = '14.84'
    D27 = '4.29'
    D28 = '11.40'
    D29 = '10.98'
    D30 = '5.94'
    D31 = '37.54'
    D32 = '7.89'
    D34 = '13.72'
    D35 = '695.79'
    D36 = '49.00'
    D37 = '67.98'
    D38 = '0.453'
    D42 = '3.40'
    D43 = '11.31'
    D44 = '0.999'
    D45 = '0.0499'
    D46 = '0.999'
    D47 = '0.153'
    D49 = '3.60'
    D50 = '3.15'
    D51 = '3.34'
}
foreach ($addr in $textForcedUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$addr]
    $cell.Style = "Normal"
}
